$d = $word.ActiveDocument

$d.Content.Find.Execute("276×7=1932", $true, $false, $false, $false, $false, $true, 1, $false, "971×6=5826", 2) | Out-Null
$d.Content.Find.Execute("856×2=1712", $true, $false, $false, $false, $false, $true, 1, $false, "936×6=5616", 2) | Out-Null
$d.Content.Find.Execute("859×5=4295", $true, $false, $false, $false, $false, $true, 1, $false, "869×6=5214", 2) | Out-Null
$d.Content.Find.Execute("437×8=3496", $true, $false, $false, $false, $false, $true, 1, $false, "905×5=4525", 2) | Out-Null
$d.Content.Find.Execute("589×9=5301", $true, $false, $false, $false, $false, $true, 1, $false, "725×7=5075", 2) | Out-Null
$d.Content.Find.Execute("181×4=724", $true, $false, $false, $false, $false, $true, 1, $false, "901×6=5406", 2) | Out-Null
$d.Content.Find.Execute("335×3=1005", $true, $false, $false, $false, $false, $true, 1, $false, "790×9=7110", 2) | Out-Null
$d.Content.Find.Execute("485×5=2425", $true, $false, $false, $false, $false, $true, 1, $false, "784×3=2352", 2) | Out-Null
$d.Content.Find.Execute("516×8=4128", $true, $false, $false, $false, $false, $true, 1, $false, "870×5=4350", 2) | Out-Null
$d.Content.Find.Execute("473×5=2365", $true, $false, $false, $false, $false, $true, 1, $false, "744×2=1488", 2) | Out-Null
$d.Content.Find.Execute("798×4=3192", $true, $false, $false, $false, $false, $true, 1, $false, "756×9=6804", 2) | Out-Null
$d.Content.Find.Execute("191×7=1337", $true, $false, $false, $false, $false, $true, 1, $false, "267×4=1068", 2) | Out-Null
$d.Content.Find.Execute("739×6=4434", $true, $false, $false, $false, $false, $true, 1, $false, "456×5=2280", 2) | Out-Null
$d.Content.Find.Execute("882×9=7938", $true, $false, $false, $false, $false, $true, 1, $false, "639×3=1917", 2) | Out-Null
$d.Content.Find.Execute("233×3=699", $true, $false, $false, $false, $false, $true, 1, $false, "633×4=2532", 2) | Out-Null
$d.Content.Find.Execute("792×9=7128", $true, $false, $false, $false, $false, $true, 1, $false, "933×7=6531", 2) | Out-Null
$d.Content.Find.Execute("297×2=594", $true, $false, $false, $false, $false, $true, 1, $false, "754×2=1508", 2) | Out-Null
$d.Content.Find.Execute("858×5=4290", $true, $false, $false, $false, $false, $true, 1, $false, "487×5=2435", 2) | Out-Null
$d.Content.Find.Execute("564×3=1692", $true, $false, $false, $false, $false, $true, 1, $false, "810×4=3240", 2) | Out-Null
$d.Content.Find.Execute("246×3=738", $true, $false, $false, $false, $false, $true, 1, $false, "419×8=3352", 2) | Out-Null
$d.Content.Find.Execute("908×5=4540", $true, $false, $false, $false, $false, $true, 1, $false, "164×3=492", 2) | Out-Null
$d.Content.Find.Execute("320×9=2880", $true, $false, $false, $false, $false, $true, 1, $false, "766×5=3830", 2) | Out-Null
$d.Content.Find.Execute("285×6=1710", $true, $false, $false, $false, $false, $true, 1, $false, "839×8=6712", 2) | Out-Null
$d.Content.Find.Execute("984×4=3936", $true, $false, $false, $false, $false, $true, 1, $false, "918×3=2754", 2) | Out-Null
$d.Content.Find.Execute("882×2=1764", $true, $false, $false, $false, $false, $true, 1, $false, "602×6=3612", 2) | Out-Null
